$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 33: effort re-booked from "Effort" + "Additional Effort" (2.75 + 2)
#     onto a single "Effort" value of 4.75; the Additional Effort cell is cleared.
$ws.Range("B33").Value = 4.75
$ws.Range("C33").Value = ""

# --- New row 58: new log entry
$ws.Range("A58").Value = "12/13/2012"
$ws.Range("B58").Value = 0.5
$ws.Range("D58").Value = "Documentation of code slightly improved"

# --- Update the view: scroll so row 17 is at the top and select C33
$win = $excel.ActiveWindow
$win.ScrollRow = 17
$win.ScrollColumn = 1
$ws.Range("C33").Select()
